$wb = $excel.ActiveWorkbook

# --- Sheet 1: _set_CASES -> _set_CASE -------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "_set_CASE"
$ws1.Cells.Item(1,1).Value = "case_Name"
# A2 already holds "baseline" - value unchanged

# --- Sheet 2: _set_YEARS ---------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(1,1).Value = "years_Name"
$ws2.Cells.Item(1,2).Value = "years_Category_1"
$ws2.Cells.Item(1,2).Copy()
$ws2.Cells.Item(1,3).PasteSpecial(-4122)
$ws2.Cells.Item(1,3).Value = "years_Aggregation"

# --- Sheet 3: _set_TECHNOLOGIES --------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(1,1).Value = "technologies_Name"
$ws3.Cells.Item(1,2).Value = "technologies_Category_1"
$ws3.Cells.Item(1,3).Value = "technologies_Category_2"
$ws3.Cells.Item(1,3).Copy()
$ws3.Cells.Item(1,4).PasteSpecial(-4122)
$ws3.Cells.Item(1,4).Value = "technologies_Aggregation"

# --- Sheet 4: _set_FLOWS -----------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(1,1).Value = "flows_Name"
$ws4.Cells.Item(1,2).Value = "flows_Category_1"
$ws4.Cells.Item(1,3).Value = "flows_Aggregation"

# --- Sheet 5: _set_FLOWS_AGG --------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Cells.Item(1,1).Value = "flows_agg_Name"
$ws5.Cells.Item(1,2).Value = "flows_agg_Category_1"
$ws5.Cells.Item(1,2).Copy()
$ws5.Cells.Item(1,3).PasteSpecial(-4122)
$ws5.Cells.Item(1,3).Value = "flows_agg_Aggregation"

# --- Selections (order matters: last Select() wins the active/tabSelected sheet) ---
$ws1.Range("A2").Select()
$ws2.Range("A2:B11").Select()
$ws3.Range("A2:C6").Select()
$ws4.Range("A2:C4").Select()
$ws5.Range("F22").Select()
